$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Insert a new column I ("category"), pushing the old I/J/K (date,
#     legislator_name, legislator_id) one column to the right. ---
$ws.Columns("I:I").Insert()

# --- Header row ---
$ws.Cells.Item(1, 9).Value = "category"

$headerFmt = $ws.Cells.Item(1, 2)

$dst = $ws.Cells.Item(1, 13)
$headerFmt.Copy($dst)
$dst.Value = "source_file"

$dst = $ws.Cells.Item(1, 14)
$headerFmt.Copy($dst)
$dst.Value = "index"

# --- Data rows: 2..26 ---
for ($r = 2; $r -le 26; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmpe4561"
    $ws.Cells.Item($r, 14).Value = $idx
}

# --- Data corrections (typo / mangled-character fixes) ---
$ws.Cells.Item(8, 6).Value = "新幣"        # F8: 新_幣 -> 新幣
$ws.Cells.Item(14, 5).Value = 10           # E14: 0.1  -> 10
$ws.Cells.Item(19, 4).Value = "1800000"    # D19: .1800000 -> 1800000
$ws.Cells.Item(21, 7).Value = "29000000"   # G21: 29000'000 -> 29000000
